# Applies the "456a3b4" data refresh to the 展览 (sheet 1) and 全部类型 (sheet 4)
# worksheets of the 江西-漫展信息 workbook:
#   - a handful of "想去人数" (F column) counters tick up
#   - a brand-new event ("南昌·花绒万兽秋镜派对", 2024-11-02) is inserted as a
#     new row right before the "上饶·星河城市动漫文化节" row, pushing every
#     row below it down by one
#   - the "南昌·云芽动漫音乐嘉年华" row's F value lands on 2254 after the bump

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Writing a literal like "2024-11-02" through .Value normally gets
    # reinterpreted by Excel as a date serial. Forcing a text number format
    # before the assignment keeps it as a plain string; clearing formats
    # afterwards drops the now-unneeded "@" format so the cell is left with
    # no explicit style, matching a freshly authored text cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Update-MandisplaySheet($ws) {
    # --- straightforward counter bumps (rows above the insertion point, so
    #     row numbers do not move) ---
    $ws.Range("F4").Value = 64
    $ws.Range("F5").Value = 1012
    $ws.Range("F6").Value = 75
    $ws.Range("F10").Value = 4713
    $ws.Range("F20").Value = 3631
    $ws.Range("F21").Value = 331
    $ws.Range("F30").Value = 218

    # --- insert the new row for 南昌·花绒万兽秋镜派对 right above the first
    #     existing row 31 entry, pushing the rest of the table down ---
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $ws.Rows.Item(31).Insert()

    # Restore the index column's style on the freshly inserted row: Insert()
    # copies formatting from the row above but drops the border, so it picks
    # up a new style slot instead of reusing the shared one every other "A"
    # cell uses. Re-applying the same formatting lets the engine fold it back
    # onto the existing style instead of leaving an orphan.
    $idxCell = $ws.Range("A31")
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    # Renumber the 0-based index column from the insertion point down to the
    # (new) last row, since column A is a running count, not copied data.
    $newLastRow = $lastRow + 1
    for ($r = 31; $r -le $newLastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Populate the new row's data.
    Set-TextValue $ws.Range("B31") "2024-11-02"
    $ws.Range("C31").Value = "南昌·花绒万兽秋镜派对"
    $ws.Range("D31").Value = "双港西大街899号 旭辉Cmall(南昌店)"
    $ws.Range("E31").Value = "2024.11.02 10:00-11.03 21:30"
    $ws.Range("F31").Value = 0
    $ws.Range("G31").Value = 168
    $ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=92859"
    $ws.Range("I31").Value = "//i2.hdslb.com/bfs/openplatform/202409/7hJL2m3F1727175584690.jpeg"

    # --- the 南昌·云芽动漫音乐嘉年华 row's F value bumps to 2254, wherever it
    #     landed after the shift (search by the row's Link id so it's robust
    #     to the exact row number on either sheet) ---
    for ($r = 31; $r -le $newLastRow; $r++) {
        $link = $ws.Cells.Item($r, 8).Value()
        if ($link -eq "https://show.bilibili.com/platform/detail.html?id=92144") {
            $ws.Cells.Item($r, 6).Value = 2254
        }
    }
}

Update-MandisplaySheet $wb.Worksheets.Item(1)
Update-MandisplaySheet $wb.Worksheets.Item(4)
